$d = $word.ActiveDocument

$pairs = @(
    @("48+9=57", "41-35=6"),
    @("34+57=91", "40-27=13"),
    @("91-43=48", "55+18=73"),
    @("40-22=18", "72-63=9"),
    @("44-7=37", "6+5=11"),
    @("16+37=53", "17-8=9"),
    @("60-1=59", "80-66=14"),
    @("5+39=44", "78+18=96"),
    @("45+46=91", "31-6=25"),
    @("19+66=85", "33-19=14"),
    @("96-39=57", "25-18=7"),
    @("29+63=92", "17+44=61"),
    @("72-57=15", "70-7=63"),
    @("53+38=91", "90-26=64"),
    @("64-17=47", "54+19=73"),
    @("49+29=78", "28+23=51"),
    @("33+29=62", "43-37=6"),
    @("5+68=73", "46-9=37"),
    @("84-29=55", "58+13=71"),
    @("69+6=75", "42-28=14"),
    @("76-57=19", "9+13=22"),
    @("18+56=74", "47+6=53"),
    @("82-27=55", "92-59=33"),
    @("29+68=97", "17+35=52"),
    @("9+14=23", "9+55=64"),
    @("17+26=43", "31-9=22"),
    @("62-8=54", "71-42=29"),
    @("40-36=4", "53-34=19"),
    @("18+44=62", "46+37=83"),
    @("19+36=55", "8+36=44"),
    @("17+46=63", "5+9=14"),
    @("57+25=82", "39+56=95"),
    @("91-75=16", "42-3=39"),
    @("63-47=16", "66+26=92"),
    @("11-5=6", "5+19=24"),
    @("52+39=91", "85-78=7"),
    @("39+26=65", "16+55=71"),
    @("52-34=18", "43+49=92"),
    @("82-75=7", "16+29=45"),
    @("29+19=48", "7+48=55"),
    @("63+18=81", "88-69=19"),
    @("70-16=54", "27+17=44"),
    @("32-24=8", "37+38=75"),
    @("46+27=73", "90-4=86"),
    @("86-38=48", "91-76=15"),
    @("46+35=81", "97-59=38"),
    @("93-65=28", "62-34=28"),
    @("50-7=43", "84-49=35"),
    @("48+39=87", "61-58=3"),
    @("38+46=84", "53-44=9"),
    @("16+77=93", "80-35=45"),
    @("7+26=33", "17+17=34"),
    @("45-39=6", "71-62=9"),
    @("16+65=81", "53+29=82"),
    @("88-39=49", "18+58=76"),
    @("37-18=19", "15+78=93"),
    @("71-38=33", "29+4=33"),
    @("65-47=18", "76+16=92"),
    @("91-13=78", "91-19=72"),
    @("74-26=48", "30-4=26"),
    @("6+6=12", "38+33=71"),
    @("83-34=49", "35+59=94"),
    @("78+8=86", "14-8=6"),
    @("57-48=9", "16+55=71"),
    @("19+28=47", "32-3=29"),
    @("62-27=35", "26+25=51"),
    @("80-77=3", "84-36=48"),
    @("44+27=71", "49+48=97"),
    @("32+9=41", "81-52=29"),
    @("39+2=41", "94-49=45"),
    @("16+75=91", "12+79=91"),
    @("35+17=52", "85-47=38"),
    @("60-38=22", "76+16=92"),
    @("2+69=71", "5+38=43"),
    @("87+7=94", "62-6=56"),
    @("9+12=21", "76-68=8"),
    @("57-19=38", "9+15=24"),
    @("38+58=96", "64-57=7"),
    @("18+15=33", "81-8=73"),
    @("35+48=83", "93-25=68"),
    @("29+55=84", "80-58=22"),
    @("45+39=84", "51-4=47"),
    @("82-57=25", "61-44=17"),
    @("14+57=71", "9+13=22"),
    @("83-18=65", "95-88=7"),
    @("71-57=14", "42-19=23"),
    @("32+49=81", "81-6=75"),
    @("76+6=82", "75-18=57"),
    @("94-28=66", "40-15=25"),
    @("67-18=49", "56+37=93"),
    @("72-44=28", "66-59=7"),
    @("29+26=55", "7+66=73"),
    @("46+25=71", "67+17=84"),
    @("73-38=35", "70-11=59"),
    @("27+67=94", "40-32=8"),
    @("55+37=92", "62-57=5"),
    @("9+34=43", "35+48=83"),
    @("96-78=18", "7+46=53"),
    @("65+16=81", "42-14=28"),
    @("63-54=9", "85-17=68")
)

foreach ($p in $pairs) {
    $old = $p[0]
    $new = $p[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Replaced $($pairs.Count) math problems"
